# Add two new columns (I: I0, J: IF) to the worksheet, mirroring the
# style/content pattern of the existing H column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy the formatting (bold/centered/bordered) from the
# existing "IP" header (H1) onto the two new header cells, then set text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-38
$i0 = @(4,6,9,3,5,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,7,4,5,4,6,6,3,6,9,8,5,3,3,4,2)
$if = @(7,6,9,7,8,6,6,6,7,6,6,3,4,6,6,6,6,6,6,5,5,4,8,6,6,7,8,8,5,8,9,8,7,6,5,5,2)

for ($r = 2; $r -le 38; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i0[$idx]
    $ws.Cells.Item($r, 10).Value = $if[$idx]
}
